$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J - copy formatting from H1 (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-34
$data = @{
    2  = @(7, 8)
    3  = @(4, 7)
    4  = @(7, 7)
    5  = @(5, 7)
    6  = @(6, 7)
    7  = @(3, 4)
    8  = @(6, 7)
    9  = @(6, 9)
    10 = @(1, 2)
    11 = @(1, 7)
    12 = @(1, 5)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 7)
    18 = @(1, 3)
    19 = @(1, 5)
    20 = @(1, 4)
    21 = @(1, 5)
    22 = @(1, 6)
    23 = @(1, 7)
    24 = @(1, 5)
    25 = @(1, 6)
    26 = @(1, 5)
    27 = @(1, 6)
    28 = @(1, 5)
    29 = @(1, 6)
    30 = @(1, 5)
    31 = @(1, 4)
    32 = @(1, 2)
    33 = @(3, 4)
    34 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
